$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.092.24'
$ws.Range("E2").Value = '  +1.30%  '

$ws.Range("D3").Value = '3.749.12'
$ws.Range("E3").Value = '  +1.24%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.40'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.08'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.72%  '

$ws.Range("D7").Value = '3.747.66'
$ws.Range("E7").Value = '  +1.36%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  +1.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.167'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.45'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.40%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.13'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.20%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000248'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.91%  '

$ws.Range("D15").Value = '4.378.21'
$ws.Range("E15").Value = '  +1.29%  '

$ws.Range("D16").Value = '3.751.40'
$ws.Range("E16").Value = '  +1.16%  '

$ws.Range("D17").Value = '69.105.12'
$ws.Range("E17").Value = '  +1.40%  '

$ws.Range("E18").Value = '  +0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.113'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.19'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.02'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +19.75%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '493.14'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.86%  '

$ws.Range("E23").Value = '  +0.57%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000151'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +6.59%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.86'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.31'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.33'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.17'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.97%  '

$ws.Range("E29").Value = '  +0.02%  '

$ws.Range("E30").Value = '  +2.83%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.49'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +4.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.05'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.52%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.58'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.64%  '

$ws.Range("D34").Value = '3.894.97'
$ws.Range("E34").Value = '  +1.28%  '

$ws.Range("E35").Value = '  +0.75%  '

$ws.Range("D36").Value = '3.685.28'
$ws.Range("E36").Value = '  +1.03%  '

$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.01'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.49%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.88'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.63%  '

$ws.Range("E40").Value = '  +1.98%  '

$ws.Range("E41").Value = '  +0.83%  '

$ws.Range("E42").Value = '  +5.20%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '432.50'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.38%  '

$ws.Range("E44").Value = '  +1.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.57'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.52'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.49%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.72'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.32'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.05%  '

$ws.Range("D50").Value = '2.792.25'
$ws.Range("E50").Value = '  +1.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0352'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.83%  '
